# Edit: remove/replace the placeholder 17th row content with two new
# Development Log entries (Driver/Observer pairing) and adjust the
# selected/view range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 ---
# Date (formatted like existing B8:B16 cells, builtin numFmtId 14)
$ws.Range("B17").Value2 = 45370
$ws.Range("B17").NumberFormat = "m/d/yy"

# Time (formatted like existing C8:C16 cells, builtin numFmtId 20)
$ws.Range("C17").Value2 = 0.41666666666666669
$ws.Range("C17").NumberFormat = "h:mm"

# Duration (hours)
$ws.Range("D17").Value2 = 1

# Who did what
$ws.Range("E17").Value = "Driver"
$ws.Range("F17").Value = "Observer"

# --- Row 18 ---
$ws.Range("B18").Value2 = 45370
$ws.Range("B18").NumberFormat = "m/d/yy"

$ws.Range("C18").Value2 = 0.5
$ws.Range("C18").NumberFormat = "h:mm"

$ws.Range("D18").Value2 = 1

$ws.Range("E18").Value = "Observer"
$ws.Range("F18").Value = "Driver"

# --- Update view state to match the post-edit selection/scroll ---
$ws.Activate()
$ws.Range("D17").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
